# The workbook originally has a single sheet named "Josh" that should be
# renamed to "Table" (the rest of the workbook's data/content is unchanged).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

if ($ws.Name -eq "Josh") {
    $ws.Name = "Table"
} elseif ($ws.Name -ne "Table") {
    # Fallback in case the active sheet isn't the one we expect - look it up
    # by its known original name instead.
    $target = $wb.Worksheets.Item("Josh")
    $target.Name = "Table"
}
